$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 531, shifting existing rows 531:548 down to 534:551
$ws.Rows("531:533").Insert()

# New record block (weekly update) - date 2023-05-29 (Excel serial 45075)
$newDate = 45075

$rowsData = @(
    @{ Row = 531; L = "Especial"; M = 16; N = 350000; O = 360000; P = 355000; S = 789 },
    @{ Row = 532; L = "Primera";  M = 20; N = 310000; O = 320000; P = 315000; S = 700 },
    @{ Row = 533; L = "Segunda";  M = 20; N = 280000; O = 290000; P = 285000; S = 633 }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value() = 8
    $ws.Cells.Item($row, 2).Value() = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value() = "Coquimbo"
    $ws.Cells.Item($row, 4).Value() = $newDate
    $ws.Cells.Item($row, 5).Value() = 4
    $ws.Cells.Item($row, 6).Value() = "Fruta"
    $ws.Cells.Item($row, 7).Value() = 100101
    $ws.Cells.Item($row, 8).Value() = "Berries"
    $ws.Cells.Item($row, 9).Value() = 100101007
    $ws.Cells.Item($row, 10).Value() = "Kiwi"
    $ws.Cells.Item($row, 11).Value() = "Hayward"
    $ws.Cells.Item($row, 12).Value() = $r.L
    $ws.Cells.Item($row, 13).Value() = $r.M
    $ws.Cells.Item($row, 14).Value() = $r.N
    $ws.Cells.Item($row, 15).Value() = $r.O
    $ws.Cells.Item($row, 16).Value() = $r.P
    $ws.Cells.Item($row, 17).Value() = "`$/bins (450 kilos)"
    $ws.Cells.Item($row, 18).Value() = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value() = $r.S
    $ws.Cells.Item($row, 20).Value() = 450
}
